$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "pion4Tests"
